$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells remain text (they look numeric, e.g. "227.53" or
# thousand-dot formatted "37.769.93") so Excel does not silently coerce them to numbers.
$priceCells = @("D2","D3","D5","D7","D9","D10","D12","D13","D14","D15","D16","D17","D18","D19","D20","D22","D26","D27","D34","D36","D40","D41","D42","D45","D46","D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.769.93"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.044.91"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "227.53"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "60.16"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("D10").Value = "0.0835"
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "2.347.45"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "14.34"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "21.45"
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").Value = "5.46"
$ws.Range("E15").Value = "  +5.48%  "
$ws.Range("D16").Value = "0.762"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "2.044.91"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "37.722.54"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "69.36"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "5.92"
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "222.74"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("E25").Value = "  +2.87%  "
$ws.Range("D26").Value = "168.95"
$ws.Range("E26").Value = "  +2.26%  "
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("E32").Value = "  +8.66%  "
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("D34").Value = "0.0602"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").Value = "6.53"
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("E37").Value = "  +4.28%  "
$ws.Range("E38").Value = "  +7.19%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "18.01"
$ws.Range("E40").Value = "  +7.39%  "
$ws.Range("D41").Value = "1.531.81"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "97.68"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D45").Value = "4.18"
$ws.Range("E45").Value = "  +7.79%  "
$ws.Range("D46").Value = "0.0894"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D51").Value = "2.236.12"
$ws.Range("E51").Value = "  +0.98%  "
